# Update the Rspo3-Lrp6 NATMI LR-pairs sheet with refreshed TPM-based values.
# New sending cluster "ECs" rows are added (rows 2-7), and the original
# FAPs-sending rows are recalculated and shifted down (rows 8-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Rspo3'
$ws.Cells.Item(2, 3).Value = 'Lrp6'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01818866666666667
$ws.Cells.Item(2, 8).Value = 0.054566
$ws.Cells.Item(2, 9).Value = 0.006403810693375696
$ws.Cells.Item(2, 10).Value = 0.006403810693375696
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 12.80871533333333
$ws.Cells.Item(2, 14).Value = 38.426146
$ws.Cells.Item(2, 15).Value = 0.1716721220213608
$ws.Cells.Item(2, 16).Value = 0.1716721220213608
$ws.Cells.Item(2, 17).Value = 0.2329734536262223
$ws.Cells.Item(2, 18).Value = 2.096761082636
$ws.Cells.Item(2, 19).Value = 0.001099355770754888
$ws.Cells.Item(2, 20).Value = 0.001099355770754888

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Rspo3'
$ws.Cells.Item(3, 3).Value = 'Lrp6'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01818866666666667
$ws.Cells.Item(3, 8).Value = 0.054566
$ws.Cells.Item(3, 9).Value = 0.006403810693375696
$ws.Cells.Item(3, 10).Value = 0.006403810693375696
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 21.85073566666667
$ws.Cells.Item(3, 14).Value = 65.552207
$ws.Cells.Item(3, 15).Value = 0.2928601395225403
$ws.Cells.Item(3, 16).Value = 0.2928601395225403
$ws.Cells.Item(3, 17).Value = 0.3974357474624445
$ws.Cells.Item(3, 18).Value = 3.576921727162
$ws.Cells.Item(3, 19).Value = 0.001875420893137942
$ws.Cells.Item(3, 20).Value = 0.001875420893137941

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Rspo3'
$ws.Cells.Item(4, 3).Value = 'Lrp6'
$ws.Cells.Item(4, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01818866666666667
$ws.Cells.Item(4, 8).Value = 0.054566
$ws.Cells.Item(4, 9).Value = 0.006403810693375696
$ws.Cells.Item(4, 10).Value = 0.006403810693375696
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.846871
$ws.Cells.Item(4, 14).Value = 32.540613
$ws.Cells.Item(4, 15).Value = 0.1453779956383313
$ws.Cells.Item(4, 16).Value = 0.1453779956383313
$ws.Cells.Item(4, 17).Value = 0.1972901209953334
$ws.Cells.Item(4, 18).Value = 1.775611088958
$ws.Cells.Item(4, 19).Value = 0.0009309731630502712
$ws.Cells.Item(4, 20).Value = 0.000930973163050271

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Rspo3'
$ws.Cells.Item(5, 3).Value = 'Lrp6'
$ws.Cells.Item(5, 4).Value = 'MuSCs'
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.01818866666666667
$ws.Cells.Item(5, 8).Value = 0.054566
$ws.Cells.Item(5, 9).Value = 0.006403810693375696
$ws.Cells.Item(5, 10).Value = 0.006403810693375696
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.389532
$ws.Cells.Item(5, 14).Value = 22.168596
$ws.Cells.Item(5, 15).Value = 0.09904011496636306
$ws.Cells.Item(5, 16).Value = 0.09904011496636304
$ws.Cells.Item(5, 17).Value = 0.1344057343706667
$ws.Cells.Item(5, 18).Value = 1.209651609336
$ws.Cells.Item(5, 19).Value = 0.000634234147294754
$ws.Cells.Item(5, 20).Value = 0.0006342341472947539

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6, 1).Value = 'ECs'
$ws.Cells.Item(6, 2).Value = 'Rspo3'
$ws.Cells.Item(6, 3).Value = 'Lrp6'
$ws.Cells.Item(6, 4).Value = 'Neutrophils'
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.01818866666666667
$ws.Cells.Item(6, 8).Value = 0.054566
$ws.Cells.Item(6, 9).Value = 0.006403810693375696
$ws.Cells.Item(6, 10).Value = 0.006403810693375696
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.190038666666666
$ws.Cells.Item(6, 14).Value = 21.570116
$ws.Cells.Item(6, 15).Value = 0.09636635393950015
$ws.Cells.Item(6, 16).Value = 0.09636635393950013
$ws.Cells.Item(6, 17).Value = 0.1307772166284445
$ws.Cells.Item(6, 18).Value = 1.176994949656
$ws.Cells.Item(6, 19).Value = 0.0006171118878393981
$ws.Cells.Item(6, 20).Value = 0.0006171118878393981

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = 'ECs'
$ws.Cells.Item(7, 2).Value = 'Rspo3'
$ws.Cells.Item(7, 3).Value = 'Lrp6'
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.01818866666666667
$ws.Cells.Item(7, 8).Value = 0.054566
$ws.Cells.Item(7, 9).Value = 0.006403810693375696
$ws.Cells.Item(7, 10).Value = 0.006403810693375696
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 14.525612
$ws.Cells.Item(7, 14).Value = 43.576836
$ws.Cells.Item(7, 15).Value = 0.1946832739119044
$ws.Cells.Item(7, 16).Value = 0.1946832739119044
$ws.Cells.Item(7, 17).Value = 0.2642015147973334
$ws.Cells.Item(7, 18).Value = 2.377813633176
$ws.Cells.Item(7, 19).Value = 0.001246714831298443
$ws.Cells.Item(7, 20).Value = 0.001246714831298443

# Row 8: FAPs -> ECs
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Rspo3'
$ws.Cells.Item(8, 3).Value = 'Lrp6'
$ws.Cells.Item(8, 4).Value = 'ECs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.822099333333334
$ws.Cells.Item(8, 8).Value = 8.466298
$ws.Cells.Item(8, 9).Value = 0.9935961893066243
$ws.Cells.Item(8, 10).Value = 0.9935961893066244
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 12.80871533333333
$ws.Cells.Item(8, 14).Value = 38.426146
$ws.Cells.Item(8, 15).Value = 0.1716721220213608
$ws.Cells.Item(8, 16).Value = 0.1716721220213608
$ws.Cells.Item(8, 17).Value = 36.14746700305645
$ws.Cells.Item(8, 18).Value = 325.327203027508
$ws.Cells.Item(8, 19).Value = 0.1705727662506059
$ws.Cells.Item(8, 20).Value = 0.170572766250606

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Rspo3'
$ws.Cells.Item(9, 3).Value = 'Lrp6'
$ws.Cells.Item(9, 4).Value = 'FAPs'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.822099333333334
$ws.Cells.Item(9, 8).Value = 8.466298
$ws.Cells.Item(9, 9).Value = 0.9935961893066243
$ws.Cells.Item(9, 10).Value = 0.9935961893066244
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 21.85073566666667
$ws.Cells.Item(9, 14).Value = 65.552207
$ws.Cells.Item(9, 15).Value = 0.2928601395225403
$ws.Cells.Item(9, 16).Value = 0.2928601395225403
$ws.Cells.Item(9, 17).Value = 61.66494655774289
$ws.Cells.Item(9, 18).Value = 554.984519019686
$ws.Cells.Item(9, 19).Value = 0.2909847186294023
$ws.Cells.Item(9, 20).Value = 0.2909847186294023

# Row 10: FAPs -> Inflammatory-Mac
$ws.Cells.Item(10, 1).Value = 'FAPs'
$ws.Cells.Item(10, 2).Value = 'Rspo3'
$ws.Cells.Item(10, 3).Value = 'Lrp6'
$ws.Cells.Item(10, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.822099333333334
$ws.Cells.Item(10, 8).Value = 8.466298
$ws.Cells.Item(10, 9).Value = 0.9935961893066243
$ws.Cells.Item(10, 10).Value = 0.9935961893066244
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.846871
$ws.Cells.Item(10, 14).Value = 32.540613
$ws.Cells.Item(10, 15).Value = 0.1453779956383313
$ws.Cells.Item(10, 16).Value = 0.1453779956383313
$ws.Cells.Item(10, 17).Value = 30.61094741785267
$ws.Cells.Item(10, 18).Value = 275.498526760674
$ws.Cells.Item(10, 19).Value = 0.144447022475281
$ws.Cells.Item(10, 20).Value = 0.144447022475281

# Row 11: FAPs -> MuSCs
$ws.Cells.Item(11, 1).Value = 'FAPs'
$ws.Cells.Item(11, 2).Value = 'Rspo3'
$ws.Cells.Item(11, 3).Value = 'Lrp6'
$ws.Cells.Item(11, 4).Value = 'MuSCs'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.822099333333334
$ws.Cells.Item(11, 8).Value = 8.466298
$ws.Cells.Item(11, 9).Value = 0.9935961893066243
$ws.Cells.Item(11, 10).Value = 0.9935961893066244
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 7.389532
$ws.Cells.Item(11, 14).Value = 22.168596
$ws.Cells.Item(11, 15).Value = 0.09904011496636306
$ws.Cells.Item(11, 16).Value = 0.09904011496636304
$ws.Cells.Item(11, 17).Value = 20.85399333084533
$ws.Cells.Item(11, 18).Value = 187.685939977608
$ws.Cells.Item(11, 19).Value = 0.0984058808190683
$ws.Cells.Item(11, 20).Value = 0.0984058808190683

# Row 12: FAPs -> Neutrophils
$ws.Cells.Item(12, 1).Value = 'FAPs'
$ws.Cells.Item(12, 2).Value = 'Rspo3'
$ws.Cells.Item(12, 3).Value = 'Lrp6'
$ws.Cells.Item(12, 4).Value = 'Neutrophils'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.822099333333334
$ws.Cells.Item(12, 8).Value = 8.466298
$ws.Cells.Item(12, 9).Value = 0.9935961893066243
$ws.Cells.Item(12, 10).Value = 0.9935961893066244
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 7.190038666666666
$ws.Cells.Item(12, 14).Value = 21.570116
$ws.Cells.Item(12, 15).Value = 0.09636635393950015
$ws.Cells.Item(12, 16).Value = 0.09636635393950013
$ws.Cells.Item(12, 17).Value = 20.29100332784089
$ws.Cells.Item(12, 18).Value = 182.619029950568
$ws.Cells.Item(12, 19).Value = 0.09574924205166074
$ws.Cells.Item(12, 20).Value = 0.09574924205166074

# Row 13: FAPs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = 'FAPs'
$ws.Cells.Item(13, 2).Value = 'Rspo3'
$ws.Cells.Item(13, 3).Value = 'Lrp6'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.822099333333334
$ws.Cells.Item(13, 8).Value = 8.466298
$ws.Cells.Item(13, 9).Value = 0.9935961893066243
$ws.Cells.Item(13, 10).Value = 0.9935961893066244
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 14.525612
$ws.Cells.Item(13, 14).Value = 43.576836
$ws.Cells.Item(13, 15).Value = 0.1946832739119044
$ws.Cells.Item(13, 16).Value = 0.1946832739119044
$ws.Cells.Item(13, 17).Value = 40.99271994145867
$ws.Cells.Item(13, 18).Value = 368.934479473128
$ws.Cells.Item(13, 19).Value = 0.193436559080606
$ws.Cells.Item(13, 20).Value = 0.1934365590806059

Write-Host "Updated dimension to A1:T13 with 12 data rows (ECs + FAPs sending clusters)."
